# Bulk User operations (#28)
# Replace the small 2-user demo sheet with the full bulk-user table:
#  - Email / Name / Role / University / Year / Group / Major / Department / Title
#  - 5 user rows with varying amounts of populated columns
#  - mailto: hyperlinks on every email cell in column A
#  - refreshed column widths + selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two existing hyperlinks (they point at the old demo addresses) so
# none of the stale relationships / styles leak into the rebuilt sheet.
$ws.Cells.Hyperlinks.Delete()

# ---- Header row ----------------------------------------------------------
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Role"
$ws.Range("D1").Value = "University"
$ws.Range("E1").Value = "Year"
$ws.Range("F1").Value = "Group"
$ws.Range("G1").Value = "Major"
$ws.Range("H1").Value = "Department"
$ws.Range("I1").Value = "Title"

# ---- Row 2: Doru5 / Management --------------------------------------------
$ws.Range("A2").Value = "bocaioandoru12+1@gmail.com"
$ws.Range("B2").Value = "Doru5"
$ws.Range("C2").Value = "Management"

# ---- Row 3: Doru1 / Professor ---------------------------------------------
$ws.Range("A3").Value = "bocaioandoru12+2@gmail.com"
$ws.Range("B3").Value = "Doru1"
$ws.Range("C3").Value = "Professor"
$ws.Range("D3").Value = "UTCN"
$ws.Range("H3").Value = "AC"
$ws.Range("I3").Value = "Prof. Eng."

# ---- Row 4: Doru2 / Student -------------------------------------------
$ws.Range("A4").Value = "bocaioandoru12+4@gmail.com"
$ws.Range("B4").Value = "Doru2"
$ws.Range("C4").Value = "Student"
$ws.Range("D4").Value = "UTCN"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "CTI"

# ---- Row 5: Doru3 / Campus_Student -----------------------------------
$ws.Range("A5").Value = "bocaioandoru12+3@gmail.com"
$ws.Range("B5").Value = "Doru3"
$ws.Range("C5").Value = "Campus_Student"
$ws.Range("D5").Value = "UTCN"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "AIA"

# ---- Row 6: Doru4 / Admin --------------------------------------------
$ws.Range("A6").Value = "bocaioandoru12+5@gmail.com"
$ws.Range("B6").Value = "Doru4"
$ws.Range("C6").Value = "Admin"

# ---- Row 7: leftover hyperlink-styled (empty) cell ------------------------
$ws.Range("A7").Style = "Hyperlink"

# ---- Hyperlinks on every e-mail cell in column A ---------------------------
# (re-apply the plain built-in "Hyperlink" cell style right after each Add so
# every linked cell shares the same style index instead of Add's private copy)
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:bocaioandoru12+4@gmail.com")
$ws.Range("A4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:bocaioandoru12+3@gmail.com")
$ws.Range("A5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:bocaioandoru12+5@gmail.com")
$ws.Range("A6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:bocaioandoru12+2@gmail.com")
$ws.Range("A3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:bocaioandoru12+1@gmail.com")
$ws.Range("A2").Style = "Hyperlink"

# ---- Column widths (best-effort match of the autofit widths) --------------
$ws.Columns.Item(1).ColumnWidth = 28.26
$ws.Columns.Item(3).ColumnWidth = 13.92
$ws.Columns.Item(4).ColumnWidth = 11.59
$ws.Columns.Item(5).ColumnWidth = 11.26
$ws.Columns.Item(8).ColumnWidth = 9.92
$ws.Columns.Item(10).ColumnWidth = 9.59

# ---- Selection restored to where the author left off -----------------------
$ws.Range("E9").Select()
